# New weekly price report: insert a new data row at the top of the
# data block (row 10), pushing all existing data rows down by one
# (old row 10 -> row 11, ... old row 65 -> row 66).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 10; Excel shifts rows 10:65 down to 11:66
# and the new row inherits formatting (incl. the date number format on D).
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with this week's record.
$ws.Range("A10").Value = 10
$ws.Range("B10").Value = "Vega Modelo de Temuco"
$ws.Range("C10").Value = "La Araucanía"
$ws.Range("D10").Value = 44532
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100101
$ws.Range("H10").Value = "Berries"
$ws.Range("I10").Value = 100101001
$ws.Range("J10").Value = "Arándano (blue)"
$ws.Range("K10").Value = "Sin especificar"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 700
$ws.Range("N10").Value = 2500
$ws.Range("O10").Value = 3000
$ws.Range("P10").Value = 2714
$ws.Range("Q10").Value = "$/kilo"
$ws.Range("R10").Value = "Región del Maule"
$ws.Range("S10").Value = 2714
$ws.Range("T10").Value = 1
